$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Sunderland vs Man City, 17:00) is removed entirely; the sheets used range shrinks to A1:AO6
$ws.Rows.Item(7).Delete()

# Rows 2-6 get refreshed League/Time/Teams/Odds data. Only cells whose value actually
# changes are written below (matching the diff) so untouched formatting/types (e.g. the
# text-stored Date column) are left exactly as-is.
# Row 2: Crystal Palace vs Fulham
$ws.Range('A2').Value = 'English Premier League'
$ws.Range('C2').Value = '14:30:00'
$ws.Range('D2').Value = 'Crystal Palace'
$ws.Range('E2').Value = 'Fulham'
$ws.Range('F2').Value = 1.27
$ws.Range('G2').Value = 1.29
$ws.Range('H2').Value = 30
$ws.Range('I2').Value = 34
$ws.Range('J2').Value = 5.1
$ws.Range('K2').Value = 5.3
$ws.Range('N2').Value = 2.12
$ws.Range('O2').Value = 1.87
$ws.Range('P2').Value = 1.21
$ws.Range('Q2').Value = 5.6
$ws.Range('R2').Value = 1.03
$ws.Range('S2').Value = 25
$ws.Range('T2').Value = 3.3
$ws.Range('U2').Value = 1.4
$ws.Range('V2').Value = 1.03
$ws.Range('W2').Value = 4.4
$ws.Range('AB2').Value = 2.34
$ws.Range('AC2').Value = 5.8
$ws.Range('AD2').Value = 36
$ws.Range('AE2').Value = 460
$ws.Range('AF2').Value = 5.3
$ws.Range('AG2').Value = 15
$ws.Range('AH2').Value = 90
$ws.Range('AJ2').Value = 30
$ws.Range('AK2').Value = 95
$ws.Range('AN2').Value = 200

# Row 3: Liverpool vs Leeds
$ws.Range('D3').Value = 'Liverpool'
$ws.Range('E3').Value = 'Leeds'
$ws.Range('F3').Value = 2.32
$ws.Range('G3').Value = 2.34
$ws.Range('H3').Value = 10.5
$ws.Range('I3').Value = 11
$ws.Range('J3').Value = 2.1
$ws.Range('K3').Value = 2.12
$ws.Range('L3').Value = 0
$ws.Range('M3').Value = 1.75
$ws.Range('N3').Value = 1.2
$ws.Range('O3').Value = 5.8
$ws.Range('P3').Value = 1.03
$ws.Range('Q3').Value = 28
$ws.Range('R3').Value = 1.01
$ws.Range('S3').Value = 120
$ws.Range('T3').Value = 12.5
$ws.Range('U3').Value = 1.07
$ws.Range('V3').Value = 1.1
$ws.Range('W3').Value = 1.65
$ws.Range('X3').Value = 2.4
$ws.Range('Z3').Value = 200
$ws.Range('AA3').Value = 1000
$ws.Range('AB3').Value = 3.3
$ws.Range('AC3').Value = 18
$ws.Range('AD3').Value = 310
$ws.Range('AE3').Value = 1000
$ws.Range('AG3').Value = 990
$ws.Range('AH3').Value = 1000
$ws.Range('AI3').Value = 1000
$ws.Range('AJ3').Value = 1000
$ws.Range('AK3').Value = 1000
$ws.Range('AL3').Value = 1000
$ws.Range('AM3').Value = 1000
$ws.Range('AN3').Value = 1000
$ws.Range('AO3').Value = 1000

# Row 4: Beitar Jerusalem vs Hapoel Tel Aviv
$ws.Range('A4').Value = 'Israeli Premier League'
$ws.Range('C4').Value = '15:30:00'
$ws.Range('D4').Value = 'Beitar Jerusalem'
$ws.Range('E4').Value = 'Hapoel Tel Aviv'
$ws.Range('F4').Value = 1.42
$ws.Range('G4').Value = 1.45
$ws.Range('H4').Value = 8.800000000000001
$ws.Range('I4').Value = 9.6
$ws.Range('J4').Value = 5
$ws.Range('K4').Value = 5.3
$ws.Range('L4').Value = 0
$ws.Range('M4').Value = 0
$ws.Range('N4').Value = 10
$ws.Range('O4').Value = 1.1
$ws.Range('P4').Value = 3
$ws.Range('Q4').Value = 1.46
$ws.Range('R4').Value = 1.67
$ws.Range('S4').Value = 2.42
$ws.Range('T4').Value = 1.5
$ws.Range('U4').Value = 2.76
$ws.Range('V4').Value = 1.11
$ws.Range('W4').Value = 3.2
$ws.Range('X4').Value = 1000
$ws.Range('Y4').Value = 1000
$ws.Range('Z4').Value = 1000
$ws.Range('AA4').Value = 1000
$ws.Range('AB4').Value = 11.5
$ws.Range('AC4').Value = 9.800000000000001
$ws.Range('AD4').Value = 22
$ws.Range('AE4').Value = 75
$ws.Range('AF4').Value = 8.4
$ws.Range('AG4').Value = 7.6
$ws.Range('AH4').Value = 15
$ws.Range('AI4').Value = 65
$ws.Range('AJ4').Value = 13.5
$ws.Range('AK4').Value = 13
$ws.Range('AL4').Value = 26
$ws.Range('AN4').Value = 9
$ws.Range('AO4').Value = 70

# Row 5: Brentford vs Tottenham
$ws.Range('A5').Value = 'English Premier League'
$ws.Range('C5').Value = '17:00:00'
$ws.Range('D5').Value = 'Brentford'
$ws.Range('E5').Value = 'Tottenham'
$ws.Range('F5').Value = 2.5
$ws.Range('G5').Value = 2.52
$ws.Range('H5').Value = 3.05
$ws.Range('I5').Value = 3.15
$ws.Range('J5').Value = 3.55
$ws.Range('K5').Value = 3.6
$ws.Range('L5').Value = 1.43
$ws.Range('M5').Value = 1.07
$ws.Range('N5').Value = 3.8
$ws.Range('O5').Value = 1.34
$ws.Range('P5').Value = 1.94
$ws.Range('Q5').Value = 2.04
$ws.Range('R5').Value = 1.37
$ws.Range('S5').Value = 3.6
$ws.Range('T5').Value = 1.79
$ws.Range('U5').Value = 2.2
$ws.Range('V5').Value = 1.46
$ws.Range('W5').Value = 1.66
$ws.Range('X5').Value = 14
$ws.Range('Y5').Value = 12
$ws.Range('Z5').Value = 19.5
$ws.Range('AA5').Value = 55
$ws.Range('AB5').Value = 10.5
$ws.Range('AC5').Value = 7.6
$ws.Range('AD5').Value = 13.5
$ws.Range('AE5').Value = 38
$ws.Range('AG5').Value = 11.5
$ws.Range('AH5').Value = 17.5
$ws.Range('AJ5').Value = 34
$ws.Range('AK5').Value = 25
$ws.Range('AL5').Value = 40
$ws.Range('AM5').Value = 90
$ws.Range('AN5').Value = 22

# Row 6: Sunderland vs Man City
$ws.Range('D6').Value = 'Sunderland'
$ws.Range('E6').Value = 'Man City'
$ws.Range('F6').Value = 8.199999999999999
$ws.Range('G6').Value = 8.6
$ws.Range('H6').Value = 1.45
$ws.Range('I6').Value = 1.46
$ws.Range('J6').Value = 5.2
$ws.Range('K6').Value = 5.3
$ws.Range('L6').Value = 1.34
$ws.Range('M6').Value = 1.04
$ws.Range('N6').Value = 4.8
$ws.Range('O6').Value = 1.25
$ws.Range('P6').Value = 2.32
$ws.Range('Q6').Value = 1.74
$ws.Range('R6').Value = 1.51
$ws.Range('S6').Value = 2.88
$ws.Range('T6').Value = 1.94
$ws.Range('U6').Value = 2.02
$ws.Range('V6').Value = 3.1
$ws.Range('W6').Value = 1.13
$ws.Range('X6').Value = 22
$ws.Range('Y6').Value = 8.4
$ws.Range('Z6').Value = 8.6
$ws.Range('AA6').Value = 12.5
$ws.Range('AB6').Value = 28
$ws.Range('AC6').Value = 11.5
$ws.Range('AD6').Value = 9.6
$ws.Range('AE6').Value = 14.5
$ws.Range('AF6').Value = 70
$ws.Range('AG6').Value = 29
$ws.Range('AH6').Value = 23
$ws.Range('AI6').Value = 32
$ws.Range('AJ6').Value = 220
$ws.Range('AK6').Value = 120
$ws.Range('AL6').Value = 95
$ws.Range('AM6').Value = 110
$ws.Range('AN6').Value = 130
$ws.Range('AO6').Value = 6.4

